$wb = $excel.ActiveWorkbook

# --- Fill in Lucas's missing round scores (previously blank, causing #DIV/0!) ---
$lucas = $wb.Worksheets.Item("Lucas")
$lucas.Range("B2").Value = 3
$lucas.Range("D2").Value = 5
$lucas.Range("F2").Value = 6
$lucas.Range("B3").Value = 2
$lucas.Range("D3").Value = 3
$lucas.Range("F3").Value = 0
$lucas.Range("B4").Value = 5
$lucas.Range("D4").Value = 5
$lucas.Range("F4").Value = 0

# --- Update the selected cell on each sheet to match the latest view state ---
$alex = $wb.Worksheets.Item("Alex")
$alex.Range("E36").Select()

$lucas.Range("F6").Select()

$grant = $wb.Worksheets.Item("Grant")
$grant.Range("E28").Select()

$alec = $wb.Worksheets.Item("Alec")
$alec.Range("H23").Select()

$scores = $wb.Worksheets.Item("Scores")
$scores.Range("D6").Select()
$scores.Activate()
